$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 95, shifting existing rows 95-118 down to 96-119
$ws.Rows(95).Insert()

# Populate new row 95 with its data (copy of surrounding data pattern, with updated values)
$ws.Cells.Item(95, 1).Value = 11
$ws.Cells.Item(95, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(95, 3).Value = "Bíobío"
$ws.Cells.Item(95, 4).Value = 44524
$ws.Cells.Item(95, 5).Value = 8
$ws.Cells.Item(95, 6).Value = "Fruta"
$ws.Cells.Item(95, 7).Value = 100108
$ws.Cells.Item(95, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(95, 9).Value = 100108005
$ws.Cells.Item(95, 10).Value = "Piña"
$ws.Cells.Item(95, 11).Value = "Caramelo"
$ws.Cells.Item(95, 12).Value = "Primera"
$ws.Cells.Item(95, 13).Value = 250
$ws.Cells.Item(95, 14).Value = 17000
$ws.Cells.Item(95, 15).Value = 18000
$ws.Cells.Item(95, 16).Value = 17400
$ws.Cells.Item(95, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(95, 18).Value = "Ecuador"
$ws.Cells.Item(95, 19).Value = 1243
$ws.Cells.Item(95, 20).Value = 14
